# Update "想去人数" (F column) values across sheets, per commit:
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 788
$ws.Range("F9").Value = 0
$ws.Range("F10").Value = 494
$ws.Range("F11").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("F18").Value = 185
$ws.Range("F20").Value = 375
$ws.Range("F21").Value = 0
$ws.Range("F22").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("F25").Value = 60
$ws.Range("F29").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("F32").Value = 556
$ws.Range("F33").Value = 49
$ws.Range("F34").Value = 2771
$ws.Range("F37").Value = 16
$ws.Range("F38").Value = 0
$ws.Range("F44").Value = 341
$ws.Range("F45").Value = 0
$ws.Range("F46").Value = 318
$ws.Range("F47").Value = 0

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 33

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 1397
$ws.Range("F4").Value = 19653
$ws.Range("F5").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("F7").Value = 1091
$ws.Range("F9").Value = 7441
$ws.Range("F10").Value = 494
$ws.Range("F11").Value = 728
$ws.Range("F12").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("F21").Value = 67
$ws.Range("F24").Value = 51
$ws.Range("F27").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("F29").Value = 11
$ws.Range("F30").Value = 168
$ws.Range("F34").Value = 49
$ws.Range("F36").Value = 2771
$ws.Range("F37").Value = 24
$ws.Range("F39").Value = 0
$ws.Range("F43").Value = 15
$ws.Range("F44").Value = 0
$ws.Range("F45").Value = 250
$ws.Range("F46").Value = 0
$ws.Range("F49").Value = 93
